# Generate Report for Handoff
# Refreshes the handoff status/date/priority values and flags a stale
# handback file across the three report sheets (Overview, zh-cn, de-de).

$wb = $excel.ActiveWorkbook

$newStatus      = "Ready for handoff"
$newHoDate      = "2016-10-25 03:25:12"   # Overview "Latest HO Xliff Generate Date" / de-de "Latest Handoff Datetime"
$newPriority    = "mt"
$newZhHoDate    = "2016-10-25 03:24:59"   # zh-cn "Latest Handoff Datetime"
$errorDetail    = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/79943214db68457fc74c8f31410552ffdf331c13/e2e/aa7bb647-5b5d-4a47-8350-bc446ca7a7d5.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ee6230d89ca8f4d1e1386ae823bf9adad6013c3f/e2e/aa7bb647-5b5d-4a47-8350-bc446ca7a7d5.md."

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")

$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

$overview.Range("G2").Value = $newHoDate
$overview.Range("G3").Value = $newHoDate

# column widths (E & F): 29.9777050018311 -> 17.2159881591797
$overview.Columns.Item(5).ColumnWidth = 16.3826530612244
$overview.Columns.Item(6).ColumnWidth = 16.3826530612244

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus

$zhcn.Range("E2").Value = $newPriority
$zhcn.Range("E3").Value = $newPriority

$zhcn.Range("H2").Value = $newZhHoDate
$zhcn.Range("H3").Value = $newZhHoDate

$zhcn.Range("P3").Value = $errorDetail

# column C: 29.9777050018311 -> 17.2159881591797
$zhcn.Columns.Item(3).ColumnWidth = 16.3826530612244
# column P: 13.7470531463623 -> 40
$zhcn.Columns.Item(16).ColumnWidth = 39.1666666666667

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus

$dede.Range("E2").Value = $newPriority
$dede.Range("E3").Value = $newPriority

$dede.Range("H2").Value = $newHoDate
$dede.Range("H3").Value = $newHoDate

$dede.Range("P3").Value = $errorDetail

# column C: 29.9777050018311 -> 17.2159881591797
$dede.Columns.Item(3).ColumnWidth = 16.3826530612244
# column P: 13.7470531463623 -> 40
$dede.Columns.Item(16).ColumnWidth = 39.1666666666667
